$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header text: "Volume 30   Number  12" -> "...13" ---
$ws.Range("A8").Characters(21, 2).Text = "13"

# --- Update week-covering dates in C9 ---
$ws.Range("C9").Characters(27, 9).Text = "3/27/2023"
$ws.Range("C9").Characters(47, 9).Text = "4/2/2023"

# --- Update the Crime Complaints data table (rows 16-27) ---
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = -25
$ws.Range("F16").Value = 6
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = -45.454545454545
$ws.Range("I16").Value = 30
$ws.Range("J16").Value = 30
$ws.Range("K16").Value = 0
$ws.Range("L16").Value = 100
$ws.Range("M16").Value = 11.111111111111
$ws.Range("N16").Value = -80.891719745222
$ws.Range("C17").Value = 2
$ws.Range("E17").Value = -33.333333333333
$ws.Range("G17").Value = 10
$ws.Range("H17").Value = -40
$ws.Range("I17").Value = 24
$ws.Range("J17").Value = 39
$ws.Range("K17").Value = -38.461538461538
$ws.Range("L17").Value = 26.315789473684
$ws.Range("M17").Value = 84.615384615384
$ws.Range("N17").Value = -54.716981132075
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 5
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 13
$ws.Range("G18").Value = 22
$ws.Range("H18").Value = -40.909090909090
$ws.Range("I18").Value = 51
$ws.Range("J18").Value = 68
$ws.Range("K18").Value = -25
$ws.Range("L18").Value = -7.272727272727
$ws.Range("M18").Value = 6.25
$ws.Range("N18").Value = -83.333333333333
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = 28.571428571428
$ws.Range("F19").Value = 38
$ws.Range("G19").Value = 43
$ws.Range("H19").Value = -11.627906976744
$ws.Range("I19").Value = 147
$ws.Range("J19").Value = 122
$ws.Range("K19").Value = 20.491803278688
$ws.Range("L19").Value = 53.125
$ws.Range("M19").Value = 96
$ws.Range("N19").Value = 67.045454545454
$ws.Range("D20").Value = 3
$ws.Range("D20").NumberFormat = "#,##0"
$ws.Range("E20").Value = 0
$ws.Range("E20").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F20").Value = 9
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = 80
$ws.Range("I20").Value = 36
$ws.Range("J20").Value = 24
$ws.Range("K20").Value = 50
$ws.Range("L20").Value = 111.764705882353
$ws.Range("M20").Value = 33.333333333333
$ws.Range("N20").Value = -83.486238532110
$ws.Range("C21").Value = 22
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 72
$ws.Range("G21").Value = 91
$ws.Range("H21").Value = -20.879120879120
$ws.Range("I21").Value = 288
$ws.Range("J21").Value = 283
$ws.Range("K21").Value = 1.766784452296
$ws.Range("L21").Value = 41.871921182266
$ws.Range("M21").Value = 50.785340314136
$ws.Range("N21").Value = -65.048543689320
$ws.Range("D22").Value = 1
$ws.Range("D22").NumberFormat = "#,##0"
$ws.Range("E22").Value = 0
$ws.Range("E22").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("F22").Value = 2
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = -33.333333333333
$ws.Range("I22").Value = 4
$ws.Range("J22").Value = 3
$ws.Range("K22").Value = 33.333333333333
$ws.Range("L22").Value = 100
$ws.Range("M22").Value = -20
$ws.Range("C23").Value = 2
$ws.Range("C23").NumberFormat = "#,##0"
$ws.Range("E23").Value = 100
$ws.Range("F23").Value = 3
$ws.Range("H23").Value = 50
$ws.Range("I23").Value = 9
$ws.Range("J23").Value = 5
$ws.Range("K23").Value = 80
$ws.Range("L23").Value = 80
$ws.Range("M23").Value = 50
$ws.Range("C24").Value = 13
$ws.Range("D24").Value = 14
$ws.Range("E24").Value = -7.142857142857
$ws.Range("F24").Value = 59
$ws.Range("H24").Value = -35.869565217391
$ws.Range("I24").Value = 211
$ws.Range("J24").Value = 251
$ws.Range("K24").Value = -15.936254980079
$ws.Range("L24").Value = 19.209039548022
$ws.Range("M24").Value = 78.813559322033
$ws.Range("C25").Value = 1
$ws.Range("D25").Value = 5
$ws.Range("E25").Value = -80
$ws.Range("F25").Value = 12
$ws.Range("H25").Value = -36.842105263157
$ws.Range("I25").Value = 52
$ws.Range("J25").Value = 65
$ws.Range("K25").Value = -20
$ws.Range("L25").Value = 8.333333333333
$ws.Range("M25").Value = 20.930232558139
# STR C27 -> si20
# STR D27 -> si20
# STR E27 -> si21
$ws.Range("F27").Value = 4
$ws.Range("H27").Value = 33.333333333333
$ws.Range("L27").Value = 87.5

# --- Row 27 (Misd. Assault): numbers -> "n/a" style text placeholders ---
# Donor cells C14/D14 (text "0", style 14) and N22 (text "***.*", style 14)
# already carry the exact formatting used for these placeholder cells.
$ws.Range("C14").Copy($ws.Range("C27"))
$ws.Range("C14").Copy($ws.Range("D27"))
$ws.Range("N22").Copy($ws.Range("E27"))
